# All Result ready to start write
#
# The sheet had a stray leading column A (values 11 / 17) that duplicated
# what ended up in column F (GENE). Removing that column shifts B:F left
# into A:E, which is exactly the new A1:E3 layout. Deleting the column
# (rather than clearing+rewriting) carries the original number formats/
# styles along for free instead of re-deriving them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()

# Fix the header typo: "MODEL_CONDITION" -> "MODELCONDITION" (now in D1
# after the column shift).
$ws.Range("D1").Value2 = "MODELCONDITION"
